$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Remove the three extra empty filler rows at the bottom (rows 17-19 before
# the edit collapse into a single trailing blank row, shifting the used
# range from A1:C19 down to A1:C16).
$ws.Rows.Item(17).Delete() | Out-Null
$ws.Rows.Item(17).Delete() | Out-Null
$ws.Rows.Item(17).Delete() | Out-Null

# Update the "Type" column wording: several rows simplify from the old
# "Opdatering og aflæsning" / "Aflæsning og beregning" / "Opdatering,
# aflæsning og beregning" phrasing down to either "Opdatering" or the new
# single word "Beregning".
$ws.Range("C4").Value = "Beregning"
$ws.Range("C5").Value = "Beregning"
$ws.Range("C7").Value = "Opdatering"
$ws.Range("C8").Value = "Opdatering"
$ws.Range("C9").Value = "Opdatering"
$ws.Range("C10").Value = "Opdatering"
$ws.Range("C11").Value = "Beregning"
$ws.Range("C12").Value = "Beregning"
$ws.Range("C13").Value = "Opdatering"
$ws.Range("C14").Value = "Beregning"
$ws.Range("C15").Value = "Beregning"

# Rows 11-15: the underlying functions/complexity were reworked -
# replace the old print/email/export/const rows with the new ones.
$ws.Range("A11").Value = "Email faktura/ordre/kontoudtog/provisionsseddel"

$ws.Range("A12").Value = "Eksporter kommaseparerede filer"

$ws.Range("A13").Value = "Opret, rediger, slet bruger"
$ws.Range("B13").Value = "Simpel"

$ws.Range("A14").Value = "Opret PDF af ordre/faktura/kontoudtog/provisionsseddel"
$ws.Range("B14").Value = "Kompleks"

$ws.Range("A15").Value = "Udregn priser på ordrer(total, total+moms)"
$ws.Range("B15").Value = "Kompleks"

# Move the active selection to A6, matching the author's cursor position
# when they saved the file.
$ws.Range("A6").Select() | Out-Null
